$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 6000
$ws.Range("J16").Value = 6000
$ws.Range("L16").Value = 6000
$ws.Range("N16").Value = -6460
$ws.Range("H98").Value = 2267.2856
$ws.Range("I98").Value = 977.5
$ws.Range("J98").Value = 10006
$ws.Range("K98").Value = 977.5
$ws.Range("L98").Value = 10006
$ws.Range("M98").Value = 520.5
$ws.Range("N98").Value = -13002
$ws.Range("H103").Value = 3525.3635
$ws.Range("I103").Value = 2857.1428
$ws.Range("J103").Value = 4694.75
$ws.Range("K103").Value = 8571.428400000001
$ws.Range("L103").Value = 14084.25
$ws.Range("M103").Value = -7985.428400000001
$ws.Range("N103").Value = -15256.25
$ws.Range("H107").Value = 160.23077
$ws.Range("I107").Value = 163.16667
$ws.Range("J107").Value = 125
$ws.Range("K107").Value = 163.16667
$ws.Range("L107").Value = 125
$ws.Range("M107").Value = 1756.83333
$ws.Range("N107").Value = -3965
$ws.Range("H122").Value = 2267.2856
$ws.Range("I122").Value = 977.5
$ws.Range("J122").Value = 10006
$ws.Range("K122").Value = 2932.5
$ws.Range("L122").Value = 30018
$ws.Range("M122").Value = -482.5
$ws.Range("N122").Value = -34918
$ws.Range("H137").Value = 1767.625
$ws.Range("I137").Value = 1598.6154
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 4795.8462
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -2245.8462
$ws.Range("N137").Value = -12600
$ws.Range("H138").Value = 3260.2222
$ws.Range("I138").Value = 8000
$ws.Range("J138").Value = 2981.4119
$ws.Range("K138").Value = 24000
$ws.Range("L138").Value = 8944.235700000001
$ws.Range("M138").Value = -18860
$ws.Range("N138").Value = -19224.2357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2877.5557
$ws.Range("I2").Value = 812.3333
$ws.Range("K2").Value = 812.3333
$ws.Range("M2").Value = -699.3333
$ws.Range("H74").Value = 745.7143
$ws.Range("I74").Value = 745.7143
$ws.Range("K74").Value = 745.7143
$ws.Range("M74").Value = 128.2857
$ws.Range("H77").Value = 745.7143
$ws.Range("I77").Value = 745.7143
$ws.Range("K77").Value = 3728.5715
$ws.Range("M77").Value = 639.4285
$ws.Range("H97").Value = 94.72727
$ws.Range("I97").Value = 111
$ws.Range("K97").Value = 111
$ws.Range("M97").Value = 385
$ws.Range("H102").Value = 2882.6365
$ws.Range("I102").Value = 1118.3334
$ws.Range("K102").Value = 1118.3334
$ws.Range("M102").Value = 503.6666
$ws.Range("H104").Value = 27000
$ws.Range("J104").Value = 27000
$ws.Range("L104").Value = 27000
$ws.Range("N104").Value = -33988
$ws.Range("H107").Value = 47499
$ws.Range("J107").Value = 47499
$ws.Range("L107").Value = 47499
$ws.Range("N107").Value = -55179
$ws.Range("H116").Value = 2877.5557
$ws.Range("I116").Value = 812.3333
$ws.Range("K116").Value = 812.3333
$ws.Range("M116").Value = 1481.6667
$ws.Range("H132").Value = 1127.6923
$ws.Range("I132").Value = 1127.6923
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3383.0769
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -853.0769
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2877.5557
$ws.Range("I3").Value = 812.3333
$ws.Range("K3").Value = 812.3333
$ws.Range("M3").Value = -698.3333
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H119").Value = 35000
$ws.Range("J119").Value = 35000
$ws.Range("L119").Value = 35000
$ws.Range("N119").Value = -44676
$ws.Range("H134").Value = 5995.8
$ws.Range("I134").Value = 5989.5
$ws.Range("K134").Value = 17968.5
$ws.Range("M134").Value = -15433.5
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -713
$ws.Range("H22").Value = 1396.8
$ws.Range("I22").Value = 859.6667
$ws.Range("J22").Value = 1627
$ws.Range("K22").Value = 859.6667
$ws.Range("L22").Value = 1627
$ws.Range("M22").Value = -509.6667
$ws.Range("N22").Value = -2327
$ws.Range("H31").Value = 2935.1667
$ws.Range("I31").Value = 2172.3
$ws.Range("K31").Value = 2172.3
$ws.Range("M31").Value = -1877.3
$ws.Range("H34").Value = 2935.1667
$ws.Range("I34").Value = 2172.3
$ws.Range("K34").Value = 2172.3
$ws.Range("M34").Value = -1970.3
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 838.4167
$ws.Range("I107").Value = 609
$ws.Range("J107").Value = 1002.2857
$ws.Range("K107").Value = 1827
$ws.Range("L107").Value = 3006.8571
$ws.Range("M107").Value = 93
$ws.Range("N107").Value = -6846.8571
$ws.Range("H128").Value = 499990
$ws.Range("I128").Value = 499990
$ws.Range("K128").Value = 1499970
$ws.Range("M128").Value = -1494990

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 44000
$ws.Range("I62").Value = 44000
$ws.Range("K62").Value = 44000
$ws.Range("M62").Value = -43314
$ws.Range("H65").Value = 44000
$ws.Range("I65").Value = 44000
$ws.Range("K65").Value = 132000
$ws.Range("M65").Value = -128568
$ws.Range("H104").Value = 23495
$ws.Range("J104").Value = 23495
$ws.Range("L104").Value = 23495
$ws.Range("N104").Value = -30483
$ws.Range("H132").Value = 3178.1
$ws.Range("I132").Value = 2575.7334
$ws.Range("K132").Value = 7727.2002
$ws.Range("M132").Value = -5197.2002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1004.63635
$ws.Range("J22").Value = 1435
$ws.Range("L22").Value = 1435
$ws.Range("N22").Value = -2025
$ws.Range("H27").Value = 1004.63635
$ws.Range("J27").Value = 1435
$ws.Range("L27").Value = 1435
$ws.Range("N27").Value = -1649
$ws.Range("H29").Value = 25007.5
$ws.Range("J29").Value = 29999
$ws.Range("L29").Value = 29999
$ws.Range("N29").Value = -30589
$ws.Range("H42").Value = 25000000
$ws.Range("I42").Value = 25000000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 25000000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -24999437
$ws.Range("N42").ClearContents()
$ws.Range("H43").Value = 10000
$ws.Range("I43").Value = 10000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 10000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -9807
$ws.Range("N43").ClearContents()
$ws.Range("H46").Value = 1111.25
$ws.Range("I46").Value = 1450
$ws.Range("J46").Value = 998.3333
$ws.Range("K46").Value = 1450
$ws.Range("L46").Value = 998.3333
$ws.Range("M46").Value = -1262
$ws.Range("N46").Value = -1374.3333
$ws.Range("H49").Value = 25000000
$ws.Range("I49").Value = 25000000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 25000000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -24999853
$ws.Range("N49").ClearContents()
$ws.Range("H55").Value = 742.2174
$ws.Range("I55").Value = 174
$ws.Range("J55").Value = 942.7646999999999
$ws.Range("K55").Value = 174
$ws.Range("L55").Value = 942.7646999999999
$ws.Range("M55").Value = -1
$ws.Range("N55").Value = -1288.7647
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H106").Value = 19288.5
$ws.Range("J106").Value = 19288.5
$ws.Range("L106").Value = 19288.5
$ws.Range("N106").Value = -21812.5
$ws.Range("H136").Value = 4046.0833
$ws.Range("I136").Value = 3965.3
$ws.Range("K136").Value = 11895.9
$ws.Range("M136").Value = -9345.900000000001
